# Listas sem duplicação de professores
# Replace cells that contained a list-like representation of teacher
# assignments (e.g. "[-, -, -, 'MEC-3B-Trat. Termicos']") with a simple "-"
# to avoid duplicated teacher entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "-"
$ws.Range("D3").Value = "-"
$ws.Range("D7").Value = "-"
$ws.Range("D8").Value = "-"
$ws.Range("C10").Value = "-"
$ws.Range("E12").Value = "-"
$ws.Range("E14").Value = "-"
$ws.Range("E15").Value = "-"
